$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 391. This shifts the existing rows
# 391..511 down to 392..512, matching the target diff exactly
# (row 391 becomes a brand new weekly record, everything below it
# keeps its previous data one row lower).
$ws.Rows("391:391").Insert()

# Populate the newly inserted row 391 with the new weekly record.
$ws.Range("A391").Value = 3
$ws.Range("B391").Value = "Femacal de La Calera"
$ws.Range("C391").Value = "Coquimbo"
$ws.Range("D391").Value = 44855
$ws.Range("E391").Value = 5
$ws.Range("F391").Value = 100112028
$ws.Range("G391").Value = "Sandia"
$ws.Range("H391").Value = "Sin especificar"
$ws.Range("I391").Value = "Primera"
$ws.Range("J391").Value = 160
$ws.Range("K391").Value = 1400
$ws.Range("L391").Value = 1400
$ws.Range("M391").Value = 1400
$ws.Range("N391").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O391").Value = "Perú"
$ws.Range("P391").Value = 1400
$ws.Range("Q391").Value = 1
$ws.Range("R391").Value = "Hortaliza"
